$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header / summary labels (sharedStrings text fixes) ---
$ws.Range("A41").Value = "Общая сумма графы “Итого”, руб."
$ws.Range("A42").Value = "Средняя площадь, кв.м."
$ws.Range("A43").Value = "Максимальный срок просрочки, дней"
$ws.Range("A44").Value = "Максимальная сумма к оплате, руб."

# --- 2. D-column tariff formula rework ---
$ws.Range("D3").Formula = "=`$A`$1 * IF(A3 <= 32, 1.1, 0.55)"
$ws.Range("D4:D38").Formula = "=`$A`$1 * IF(A4 <= 32, 1.1, 0.55)"

# --- 3. Remove the blank spacer row above the summary block, shifting the
#        summary rows (old 41-44) up to 40-43, then restore the trailing
#        blank row so the sheet keeps the same overall row count. ---
$ws.Rows(40).Delete()
$ws.Rows(61).RowHeight = 15.75

# --- 4. Move the summary labels/values one column to the right (A/B -> B/C) ---
$ws.Range("A40:B40").Cut($ws.Range("B40"))
$ws.Range("A41:B41").Cut($ws.Range("B41"))
$ws.Range("A42:B42").Cut($ws.Range("B42"))
$ws.Range("A43:B43").Cut($ws.Range("B43"))

# --- 5. View state ---
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("B44").Select()
